$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# A new row (componentId props description) is inserted right above the
# existing row 91 ("XML2SOURCE_FILE.PROPS_INTERFACE.DESCRIPTION"), pushing
# every row from 91..114 down by one (91->92 ... 114->115).
#
# We implement the insert by shifting the row contents downward (copying
# row r into row r+1, working from the bottom up so nothing is clobbered
# before it's been read), then overwriting row 91 with the new data.
# ---------------------------------------------------------------------------

for ($r = 114; $r -ge 91; $r--) {
    $src = $ws.Range("A" + $r + ":G" + $r)
    $dst = $ws.Range("A" + ($r + 1) + ":G" + ($r + 1))
    # Copy() only overwrites non-blank source cells, so clear the
    # destination first or stale values survive where the source is blank.
    $dst.ClearContents()
    $src.Copy($dst)
}

# Column A holds a shared "running number" formula (=previous+1). The plain
# Copy() above duplicates the cached value but not the relative formula, so
# restore it explicitly for every data row affected by the shift.
for ($r = 82; $r -le 113; $r++) {
    $ws.Range("A" + $r).Formula = "=A" + ($r - 1) + "+1"
}

# New row 91: componentId description key/value (new shared strings).
$ws.Range("B91").Value2 = "XML2SOURCE_FILE.PROPS_COMPOENENT_ID.DESCRIPTION"
$ws.Range("C91").Value2 = "コンポーネントのcomponentIdプロパティです"

# The sheet had scrolled to A69 with C100 selected; the saved view now shows
# the top of the sheet with B17 selected.
$ws.Range("B17").Select()
